$wb = $excel.ActiveWorkbook

# ---- PIR sheet: append rows 140-152 ----
$wsPIR = $wb.Worksheets.Item("PIR")
# Pre-format Date/Timestamp/Hour columns as Text so values like "2026-01-30"
# and "15:49:27" are stored literally instead of being auto-converted to
# date/time serial numbers.
$wsPIR.Range("A140:C152").NumberFormat = "@"

$wsPIR.Range("A140").Value = "2026-01-30"
$wsPIR.Range("B140").Value = "15:49:27"
$wsPIR.Range("C140").Value = "15:00"
$wsPIR.Range("D140").Value = "Bathroom"
$wsPIR.Range("E140").Value = "No Motion"
$wsPIR.Range("F140").Value = "Inactive"

$wsPIR.Range("A141").Value = "2026-01-30"
$wsPIR.Range("B141").Value = "15:49:30"
$wsPIR.Range("C141").Value = "15:00"
$wsPIR.Range("D141").Value = "Bathroom"
$wsPIR.Range("E141").Value = "No Motion"
$wsPIR.Range("F141").Value = "Inactive"

$wsPIR.Range("A142").Value = "2026-01-30"
$wsPIR.Range("B142").Value = "15:49:35"
$wsPIR.Range("C142").Value = "15:00"
$wsPIR.Range("D142").Value = "Bathroom"
$wsPIR.Range("E142").Value = "No Motion"
$wsPIR.Range("F142").Value = "Inactive"

$wsPIR.Range("A143").Value = "2026-01-30"
$wsPIR.Range("B143").Value = "15:49:40"
$wsPIR.Range("C143").Value = "15:00"
$wsPIR.Range("D143").Value = "Bathroom"
$wsPIR.Range("E143").Value = "No Motion"
$wsPIR.Range("F143").Value = "Inactive"

$wsPIR.Range("A144").Value = "2026-01-30"
$wsPIR.Range("B144").Value = "15:49:45"
$wsPIR.Range("C144").Value = "15:00"
$wsPIR.Range("D144").Value = "Bathroom"
$wsPIR.Range("E144").Value = "No Motion"
$wsPIR.Range("F144").Value = "Inactive"

$wsPIR.Range("A145").Value = "2026-01-30"
$wsPIR.Range("B145").Value = "15:49:50"
$wsPIR.Range("C145").Value = "15:00"
$wsPIR.Range("D145").Value = "Bathroom"
$wsPIR.Range("E145").Value = "No Motion"
$wsPIR.Range("F145").Value = "Inactive"

$wsPIR.Range("A146").Value = "2026-01-30"
$wsPIR.Range("B146").Value = "15:49:55"
$wsPIR.Range("C146").Value = "15:00"
$wsPIR.Range("D146").Value = "Bathroom"
$wsPIR.Range("E146").Value = "No Motion"
$wsPIR.Range("F146").Value = "Inactive"

$wsPIR.Range("A147").Value = "2026-01-30"
$wsPIR.Range("B147").Value = "15:50:00"
$wsPIR.Range("C147").Value = "15:00"
$wsPIR.Range("D147").Value = "Bathroom"
$wsPIR.Range("E147").Value = "No Motion"
$wsPIR.Range("F147").Value = "Inactive"

$wsPIR.Range("A148").Value = "2026-01-30"
$wsPIR.Range("B148").Value = "15:50:05"
$wsPIR.Range("C148").Value = "15:00"
$wsPIR.Range("D148").Value = "Bathroom"
$wsPIR.Range("E148").Value = "No Motion"
$wsPIR.Range("F148").Value = "Inactive"

$wsPIR.Range("A149").Value = "2026-01-30"
$wsPIR.Range("B149").Value = "15:50:10"
$wsPIR.Range("C149").Value = "15:00"
$wsPIR.Range("D149").Value = "Bathroom"
$wsPIR.Range("E149").Value = "No Motion"
$wsPIR.Range("F149").Value = "Inactive"

$wsPIR.Range("A150").Value = "2026-01-30"
$wsPIR.Range("B150").Value = "15:50:15"
$wsPIR.Range("C150").Value = "15:00"
$wsPIR.Range("D150").Value = "Bathroom"
$wsPIR.Range("E150").Value = "No Motion"
$wsPIR.Range("F150").Value = "Inactive"

$wsPIR.Range("A151").Value = "2026-01-30"
$wsPIR.Range("B151").Value = "15:50:20"
$wsPIR.Range("C151").Value = "15:00"
$wsPIR.Range("D151").Value = "Bathroom"
$wsPIR.Range("E151").Value = "No Motion"
$wsPIR.Range("F151").Value = "Inactive"

$wsPIR.Range("A152").Value = "2026-01-30"
$wsPIR.Range("B152").Value = "15:50:25"
$wsPIR.Range("C152").Value = "15:00"
$wsPIR.Range("D152").Value = "Bathroom"
$wsPIR.Range("E152").Value = "No Motion"
$wsPIR.Range("F152").Value = "Inactive"

# ---- Humidity sheet: append rows 81-88 ----
$wsHumidity = $wb.Worksheets.Item("Humidity")
# Pre-format Date/Timestamp/Hour/Value columns as Text so values like
# "2026-01-30", "15:49:28" and "87.7%" are stored literally instead of being
# auto-converted to date/time/percentage numbers.
$wsHumidity.Range("A81:C88").NumberFormat = "@"
$wsHumidity.Range("E81:E88").NumberFormat = "@"

$wsHumidity.Range("A81").Value = "2026-01-30"
$wsHumidity.Range("B81").Value = "15:49:28"
$wsHumidity.Range("C81").Value = "15:00"
$wsHumidity.Range("D81").Value = "Bathroom"
$wsHumidity.Range("E81").Value = "87.7%"
$wsHumidity.Range("F81").Value = "Active"

$wsHumidity.Range("A82").Value = "2026-01-30"
$wsHumidity.Range("B82").Value = "15:49:30"
$wsHumidity.Range("C82").Value = "15:00"
$wsHumidity.Range("D82").Value = "Bathroom"
$wsHumidity.Range("E82").Value = "87.8%"
$wsHumidity.Range("F82").Value = "Active"

$wsHumidity.Range("A83").Value = "2026-01-30"
$wsHumidity.Range("B83").Value = "15:49:35"
$wsHumidity.Range("C83").Value = "15:00"
$wsHumidity.Range("D83").Value = "Bathroom"
$wsHumidity.Range("E83").Value = "86.3%"
$wsHumidity.Range("F83").Value = "Active"

$wsHumidity.Range("A84").Value = "2026-01-30"
$wsHumidity.Range("B84").Value = "15:49:46"
$wsHumidity.Range("C84").Value = "15:00"
$wsHumidity.Range("D84").Value = "Bathroom"
$wsHumidity.Range("E84").Value = "87.7%"
$wsHumidity.Range("F84").Value = "Active"

$wsHumidity.Range("A85").Value = "2026-01-30"
$wsHumidity.Range("B85").Value = "15:49:56"
$wsHumidity.Range("C85").Value = "15:00"
$wsHumidity.Range("D85").Value = "Bathroom"
$wsHumidity.Range("E85").Value = "86.8%"
$wsHumidity.Range("F85").Value = "Active"

$wsHumidity.Range("A86").Value = "2026-01-30"
$wsHumidity.Range("B86").Value = "15:50:11"
$wsHumidity.Range("C86").Value = "15:00"
$wsHumidity.Range("D86").Value = "Bathroom"
$wsHumidity.Range("E86").Value = "87.8%"
$wsHumidity.Range("F86").Value = "Active"

$wsHumidity.Range("A87").Value = "2026-01-30"
$wsHumidity.Range("B87").Value = "15:50:16"
$wsHumidity.Range("C87").Value = "15:00"
$wsHumidity.Range("D87").Value = "Bathroom"
$wsHumidity.Range("E87").Value = "86.8%"
$wsHumidity.Range("F87").Value = "Active"

$wsHumidity.Range("A88").Value = "2026-01-30"
$wsHumidity.Range("B88").Value = "15:50:26"
$wsHumidity.Range("C88").Value = "15:00"
$wsHumidity.Range("D88").Value = "Bathroom"
$wsHumidity.Range("E88").Value = "87.7%"
$wsHumidity.Range("F88").Value = "Active"
